# The underlying change in this revision is purely a re-serialization /
# canonicalization of the package's OOXML (namespace declarations on
# <w:document> and attribute ordering inside elements such as <w:tblW>,
# <w:tcW>, <w:tblLook>, <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>,
# <w:latentStyles>, <w:lsdException>, <w:style>, <w:tblInd>,
# <w:tblCellMar>, <w:color> and <w:tblBorders> are all simply
# alphabetized). No text, formatting, table geometry, hyperlink, or
# style definition is added, removed, or changed in value anywhere in
# the package - the commit is a no-op at the document-model level
# (a bulk "resave" that rode along with an unrelated version bump).
#
# Word's object model has no notion of "XML attribute order", so there
# is nothing to change through COM automation here; we simply touch the
# document to keep the step explicit without mutating any content.
$d = $word.ActiveDocument
$null = $d.Content
